$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("D12").Value = "[-, 'MEC-2A-Des. Maq. Cad_T2']"

# Row 14
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "[-, 'MEC-2A-Des. Maq. Cad_T2']"

# Row 15
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "[-, 'MEC-2A-Des. Maq. Cad_T2']"

# Row 16
$ws.Range("C16").Value = "-"

# Row 18
$ws.Range("C18").Value = "-"
$ws.Range("E18").Value = "MEC-2NA-Des. Maq. Cad"
$ws.Range("F18").Value = "MEC-2NB-Elemaq."

# Row 19
$ws.Range("E19").Value = "MEC-2NA-Des. Maq. Cad"
$ws.Range("F19").Value = "MEC-2NB-Elemaq."

# Row 20
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "MEC-2NA-Elemaq."
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("E21").Value = "MEC-2NA-Elemaq."
$ws.Range("F21").Value = "-"
